# Apply weekly refresh of Fruta/Hortaliza prices: rotate the data in
# columns D (Fecha), M (Volumen), N (Precio mínimo), O (Precio máximo),
# P (Precio promedio ponderado), Q (Unidad de comercialización) and
# S (Precio $/Kg) across rows 2-9, as described by the commit "Fruta /
# hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row after the edit (row => values)
$data = @{
    2 = @{ D = 44232; M = 60; N = 11000; O = 12000; P = 11583; Q = "`$/caja 14 kilos empedrada"; S = 827 }
    3 = @{ D = 44172; M = 90; N = 8500;  O = 9000;  P = 8806;  Q = "`$/caja 14 kilos empedrada"; S = 629 }
    4 = @{ D = 44229; M = 55; N = 11000; O = 12000; P = 11364; Q = "`$/caja 14 kilos empedrada"; S = 812 }
    5 = @{ D = 44181; M = 65; N = 9000;  O = 10000; P = 9462;  Q = "`$/caja 14 kilos empedrada"; S = 676 }
    6 = @{ D = 44216; M = 55; N = 11000; O = 12000; P = 11545; Q = "`$/caja 14 kilos empedrada"; S = 825 }
    7 = @{ D = 44210; M = 70; N = 10000; O = 11000; P = 10357; Q = "`$/caja 14 kilos empedrada"; S = 740 }
    8 = @{ D = 45138; M = 50; N = 14000; O = 14000; P = 14000; Q = "`$/caja 14 kilos granel";    S = 1000 }
    9 = @{ D = 44253; M = 90; N = 12000; O = 13000; P = 12667; Q = "`$/caja 14 kilos empedrada"; S = 905 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]

    $ws.Cells.Item($row, 4).Value = $vals.D   # D: Fecha
    $ws.Cells.Item($row, 13).Value = $vals.M  # M: Volumen
    $ws.Cells.Item($row, 14).Value = $vals.N  # N: Precio minimo
    $ws.Cells.Item($row, 15).Value = $vals.O  # O: Precio maximo
    $ws.Cells.Item($row, 16).Value = $vals.P  # P: Precio promedio ponderado
    $ws.Cells.Item($row, 17).Value = $vals.Q  # Q: Unidad de comercializacion
    $ws.Cells.Item($row, 19).Value = $vals.S  # S: Precio $/Kg
}
